$wb = $excel.ActiveWorkbook

# --- Step 1: remove the "Signal_Value_125" column (last signal-value column) ---
# Column AL is index 38 (A=1 ... AL=38) in both Step1_Data and Step2_Sj sheets.
$ws1 = $wb.Worksheets.Item("Step1_Data")
$ws2 = $wb.Worksheets.Item("Step2_Sj")
$ws1.Columns.Item(38).Delete()
$ws2.Columns.Item(38).Delete()

# --- Step 2: refresh recomputed values for "signal segment 6" (row 7) ---
# Step1_Data (per-segment normalized signal distribution)
$ws1.Range("F7").Value = 0.005298349535039483
$ws1.Range("G7").Value = 0.1508483242301047
$ws1.Range("H7").Value = 0.1405326493782038
$ws1.Range("I7").Value = 0.3120839732408696
$ws1.Range("J7").Value = 0.00004906539731378454
$ws1.Range("K7").Value = 0.1087232357549484
$ws1.Range("L7").Value = 0.0275238799838982
$ws1.Range("M7").Value = 0.01899223325337546
$ws1.Range("N7").Value = 0.002350969720240833
$ws1.Range("O7").Value = 0.001171014874409122
$ws1.Range("P7").Value = 0.02655490458702164
$ws1.Range("Q7").Value = 0.01098517944578718
$ws1.Range("R7").Value = 0.002488528353416283
$ws1.Range("S7").Value = 0.003005630063456896
$ws1.Range("T7").Value = 0.03798160106064788
$ws1.Range("U7").Value = 0.01653792080614259
$ws1.Range("V7").Value = 0.02855530006392817
$ws1.Range("W7").Value = 0.0275238799838982
$ws1.Range("X7").Value = 0.0002251606381524403
$ws1.Range("Y7").Value = 0.0001545203383064322
$ws1.Range("Z7").Value = 0.000799344652040642
$ws1.Range("AA7").Value = 0.00924911035094059
$ws1.Range("AB7").Value = 0.0000198539713962728
$ws1.Range("AC7").Value = 0.003319000346474372
$ws1.Range("AD7").Value = 0.00257940914509669
$ws1.Range("AE7").Value = 0.0003855244476339855
$ws1.Range("AF7").Value = 0.0001582656416773186
$ws1.Range("AG7").Value = 0.005661147523478109
$ws1.Range("AH7").Value = 0.01540379514759301
$ws1.Range("AI7").Value = 0.009698609299318628
$ws1.Range("AJ7").Value = 0.003958021828773872
$ws1.Range("AK7").Value = 0.02718159693641543

# Step2_Sj (cumulative sum of Step1_Data)
$ws2.Range("F7").Value = 0.005298349535039483
$ws2.Range("G7").Value = 0.1561466737651442
$ws2.Range("H7").Value = 0.296679323143348
$ws2.Range("I7").Value = 0.6087632963842176
$ws2.Range("J7").Value = 0.6088123617815313
$ws2.Range("K7").Value = 0.7175355975364797
$ws2.Range("L7").Value = 0.7450594775203778
$ws2.Range("M7").Value = 0.7640517107737533
$ws2.Range("N7").Value = 0.7664026804939942
$ws2.Range("O7").Value = 0.7675736953684033
$ws2.Range("P7").Value = 0.7941285999554248
$ws2.Range("Q7").Value = 0.805113779401212
$ws2.Range("R7").Value = 0.8076023077546283
$ws2.Range("S7").Value = 0.8106079378180852
$ws2.Range("T7").Value = 0.848589538878733
$ws2.Range("U7").Value = 0.8651274596848756
$ws2.Range("V7").Value = 0.8936827597488037
$ws2.Range("W7").Value = 0.9212066397327019
$ws2.Range("X7").Value = 0.9214318003708544
$ws2.Range("Y7").Value = 0.9215863207091608
$ws2.Range("Z7").Value = 0.9223856653612015
$ws2.Range("AA7").Value = 0.931634775712142
$ws2.Range("AB7").Value = 0.9316546296835383
$ws2.Range("AC7").Value = 0.9349736300300127
$ws2.Range("AD7").Value = 0.9375530391751093
$ws2.Range("AE7").Value = 0.9379385636227433
$ws2.Range("AF7").Value = 0.9380968292644206
$ws2.Range("AG7").Value = 0.9437579767878987
$ws2.Range("AH7").Value = 0.9591617719354918
$ws2.Range("AI7").Value = 0.9688603812348104
$ws2.Range("AJ7").Value = 0.9728184030635844
$ws2.Range("AK7").Value = 0.9999999999999998

# --- Step 3: refresh dependent lookups in the Step3_DataPts_* sheets (row 7 = segment 6) ---
$ws3_05 = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws3_07 = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws3_08 = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws3_09 = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws3_05.Range("F7").Value = 0.6087632963842176
$ws3_07.Range("F7").Value = 0.7175355975364797
$ws3_08.Range("D7").Value = 16
$ws3_08.Range("F7").Value = 0.805113779401212
$ws3_08.Range("G7").Value = 12
$ws3_09.Range("F7").Value = 0.9212066397327019
